$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -8.01
$ws.Range("B7").Value = 5.481
$ws.Range("A8").Value = -22.366
$ws.Range("A10").Value = -21.729
$ws.Range("E10").Value = 16.414
$ws.Range("A12").Value = -21.511
$ws.Range("E12").Value = 17.858
$ws.Range("E13").Value = 16.636
$ws.Range("E14").Value = 16.805
$ws.Range("B15").Value = 5.313000000000001
$ws.Range("A18").Value = -21.808
$ws.Range("B18").Value = 5.92
$ws.Range("D18").Value = -8.753000000000002
$ws.Range("D19").Value = -7.912000000000001
$ws.Range("B20").Value = 7.489
$ws.Range("D27").Value = -7.816
$ws.Range("B29").Value = 5.377
$ws.Range("E29").Value = 16.924
$ws.Range("B30").Value = 6.000999999999999
$ws.Range("B31").Value = 5.231
$ws.Range("D31").Value = -7.606999999999999
$ws.Range("E32").Value = 16.717
$ws.Range("E35").Value = 16.194
$ws.Range("A37").Value = -20.272
$ws.Range("D38").Value = -8.286000000000001
$ws.Range("B40").Value = 8.962
$ws.Range("D42").Value = -8.257999999999999
$ws.Range("E43").Value = 16.791
$ws.Range("D44").Value = -7.936999999999999
$ws.Range("D47").Value = -7.847
$ws.Range("E48").Value = 16.954
$ws.Range("E49").Value = 16.359
$ws.Range("B50").Value = 5.012
$ws.Range("E50").Value = 16.417
$ws.Range("A55").Value = -21.843
$ws.Range("E56").Value = 16.183
$ws.Range("D58").Value = -8.371
$ws.Range("D65").Value = -7.851999999999999
$ws.Range("A68").Value = -21.454
$ws.Range("B68").Value = 5.298000000000001
$ws.Range("E69").Value = 17.267
$ws.Range("D73").Value = -8.027000000000001
$ws.Range("B76").Value = 6.33
$ws.Range("A77").Value = -20.637
$ws.Range("A78").Value = -20.126
$ws.Range("A81").Value = -21.258
$ws.Range("E81").Value = 16.673
$ws.Range("A82").Value = -22.077
$ws.Range("B87").Value = 4.636
$ws.Range("B88").Value = 4.858000000000001
$ws.Range("D90").Value = -8.211
$ws.Range("E92").Value = 17.637
$ws.Range("D94").Value = -7.486999999999999
$ws.Range("D95").Value = -7.561
$ws.Range("B96").Value = 6.281000000000001
$ws.Range("B98").Value = 5.762
$ws.Range("B101").Value = 7.875
$ws.Range("D101").Value = -8.187000000000001
$ws.Range("B102").Value = 7.375999999999999
